$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 257.66666
$ws.Range("I28").Value = 236.7
$ws.Range("K28").Value = 236.7
$ws.Range("M28").Value = 248.3
$ws.Range("H93").Value = 32243.908
$ws.Range("J93").Value = 32243.908
$ws.Range("L93").Value = 32243.908
$ws.Range("N93").Value = -37235.908
$ws.Range("H105").Value = 49671
$ws.Range("J105").Value = 49671
$ws.Range("L105").Value = 49671
$ws.Range("N105").Value = -56659
$ws.Range("H107").Value = 1388.75
$ws.Range("I107").Value = 1388.75
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1388.75
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 531.25
$ws.Range("N107").ClearContents()
$ws.Range("H126").Value = 42243.2
$ws.Range("J126").Value = 42243.2
$ws.Range("L126").Value = 42243.2
$ws.Range("N126").Value = -52123.2
$ws.Range("H128").Value = 35219.2
$ws.Range("J128").Value = 35219.2
$ws.Range("L128").Value = 35219.2
$ws.Range("N128").Value = -45179.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1235.3636
$ws.Range("I74").Value = 681.9474
$ws.Range("J74").Value = 1986.4286
$ws.Range("K74").Value = 681.9474
$ws.Range("L74").Value = 1986.4286
$ws.Range("M74").Value = 192.0526
$ws.Range("N74").Value = -3734.4286
$ws.Range("H77").Value = 1235.3636
$ws.Range("I77").Value = 681.9474
$ws.Range("J77").Value = 1986.4286
$ws.Range("K77").Value = 3409.737
$ws.Range("L77").Value = 9932.143
$ws.Range("M77").Value = 958.2629999999999
$ws.Range("N77").Value = -18668.143
$ws.Range("H101").Value = 48560
$ws.Range("J101").Value = 48560
$ws.Range("L101").Value = 48560
$ws.Range("N101").Value = -55050
$ws.Range("H102").Value = 58837576
$ws.Range("I102").Value = 71430690
$ws.Range("J102").Value = 69714
$ws.Range("K102").Value = 71430690
$ws.Range("L102").Value = 69714
$ws.Range("M102").Value = -71429068
$ws.Range("N102").Value = -72958
$ws.Range("H105").Value = 43742.8
$ws.Range("J105").Value = 43742.8
$ws.Range("L105").Value = 43742.8
$ws.Range("N105").Value = -50730.8
$ws.Range("H106").Value = 39183
$ws.Range("J106").Value = 39183
$ws.Range("L106").Value = 39183
$ws.Range("N106").Value = -41707
$ws.Range("H110").Value = 2028.9286
$ws.Range("I110").Value = 1954.5834
$ws.Range("K110").Value = 1954.5834
$ws.Range("M110").Value = 90.41660000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H113").Value = 2500
$ws.Range("I113").Value = 2500
$ws.Range("K113").Value = 2500
$ws.Range("M113").Value = -330

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 39549.332
$ws.Range("J43").Value = 39549.332
$ws.Range("L43").Value = 39549.332
$ws.Range("N43").Value = -39917.332
$ws.Range("H57").Value = 34957
$ws.Range("J57").Value = 34957
$ws.Range("L57").Value = 34957
$ws.Range("N57").Value = -36077
$ws.Range("H101").Value = 39549.332
$ws.Range("J101").Value = 39549.332
$ws.Range("L101").Value = 39549.332
$ws.Range("N101").Value = -46039.332
$ws.Range("H124").Value = 45318
$ws.Range("J124").Value = 45318
$ws.Range("L124").Value = 45318
$ws.Range("N124").Value = -50228
$ws.Range("H125").Value = 30664
$ws.Range("J125").Value = 30664
$ws.Range("L125").Value = 30664
$ws.Range("N125").Value = -35584
$ws.Range("H131").Value = 35621
$ws.Range("J131").Value = 35621
$ws.Range("L131").Value = 35621
$ws.Range("N131").Value = -45701

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1972.2727
$ws.Range("I97").Value = 1463.4706
$ws.Range("J97").Value = 3702.2
$ws.Range("K97").Value = 1463.4706
$ws.Range("L97").Value = 3702.2
$ws.Range("M97").Value = -967.4706000000001
$ws.Range("N97").Value = -4694.2
$ws.Range("H101").Value = 50649
$ws.Range("J101").Value = 50649
$ws.Range("L101").Value = 50649
$ws.Range("N101").Value = -57139
$ws.Range("H113").Value = 2622.2222
$ws.Range("I113").Value = 2240
$ws.Range("J113").Value = 3100
$ws.Range("K113").Value = 2240
$ws.Range("L113").Value = 3100
$ws.Range("M113").Value = -70
$ws.Range("N113").Value = -7440
$ws.Range("H118").Value = 38298
$ws.Range("J118").Value = 38298
$ws.Range("L118").Value = 38298
$ws.Range("N118").Value = -41612
$ws.Range("H120").Value = 39309
$ws.Range("J120").Value = 39309
$ws.Range("L120").Value = 39309
$ws.Range("N120").Value = -48985
$ws.Range("H125").Value = 40992
$ws.Range("J125").Value = 40992
$ws.Range("L125").Value = 40992
$ws.Range("N125").Value = -45912
$ws.Range("H127").Value = 41996
$ws.Range("J127").Value = 41996
$ws.Range("L127").Value = 41996
$ws.Range("N127").Value = -51916
$ws.Range("H131").Value = 38986
$ws.Range("J131").Value = 38986
$ws.Range("L131").Value = 38986
$ws.Range("N131").Value = -49066
$ws.Range("H139").Value = 21845.2
$ws.Range("J139").Value = 21845.2
$ws.Range("L139").Value = 21845.2
$ws.Range("N139").Value = -32125.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2599.75
$ws.Range("I40").Value = 2699.5
$ws.Range("J40").Value = 2500
$ws.Range("K40").Value = 2699.5
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = -2563.5
$ws.Range("N40").Value = -2772
$ws.Range("H105").Value = 34084
$ws.Range("J105").Value = 34084
$ws.Range("L105").Value = 34084
$ws.Range("N105").Value = -41072
$ws.Range("H109").Value = 35281
$ws.Range("J109").Value = 35281
$ws.Range("L109").Value = 35281
$ws.Range("N109").Value = -38055
$ws.Range("H117").Value = 45380
$ws.Range("J117").Value = 45380
$ws.Range("L117").Value = 45380
$ws.Range("N117").Value = -54558
$ws.Range("H120").Value = 42143.8
$ws.Range("J120").Value = 42143.8
$ws.Range("L120").Value = 42143.8
$ws.Range("N120").Value = -51819.8
$ws.Range("H123").Value = 32878
$ws.Range("J123").Value = 32878
$ws.Range("L123").Value = 32878
$ws.Range("N123").Value = -42678
$ws.Range("H129").Value = 32404.834
$ws.Range("J129").Value = 32404.834
$ws.Range("L129").Value = 32404.834
$ws.Range("N129").Value = -42404.834
$ws.Range("H131").Value = 43326
$ws.Range("J131").Value = 43326
$ws.Range("L131").Value = 43326
$ws.Range("N131").Value = -53406

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 27406.5
$ws.Range("J27").Value = 27406.5
$ws.Range("L27").Value = 27406.5
$ws.Range("N27").Value = -27544.5
$ws.Range("H103").Value = 39126
$ws.Range("J103").Value = 39126
$ws.Range("L103").Value = 39126
$ws.Range("N103").Value = -41470
$ws.Range("H104").Value = 39185
$ws.Range("J104").Value = 39185
$ws.Range("L104").Value = 39185
$ws.Range("N104").Value = -46173
$ws.Range("H109").Value = 30639
$ws.Range("J109").Value = 30639
$ws.Range("L109").Value = 30639
$ws.Range("N109").Value = -33413
$ws.Range("H118").Value = 26467.5
$ws.Range("J118").Value = 30290
$ws.Range("L118").Value = 30290
$ws.Range("N118").Value = -33604
$ws.Range("H123").Value = 31610.5
$ws.Range("J123").Value = 31610.5
$ws.Range("L123").Value = 31610.5
$ws.Range("N123").Value = -41410.5
$ws.Range("H127").Value = 16765.363
$ws.Range("J127").Value = 29473
$ws.Range("L127").Value = 29473
$ws.Range("N127").Value = -39393
$ws.Range("H129").Value = 33192
$ws.Range("J129").Value = 33192
$ws.Range("L129").Value = 33192
$ws.Range("N129").Value = -43192
